$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '257.40'
$c.ClearFormats()

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = '@'
$c.Value = '-0.44%'
$c.ClearFormats()

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = '@'
$c.Value = '0.22%'
$c.ClearFormats()

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '4.579'
$c.ClearFormats()

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = '@'
$c.Value = '-6.22%'
$c.ClearFormats()

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '0.05898'
$c.ClearFormats()

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = '@'
$c.Value = '-1.05%'
$c.ClearFormats()

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '6.632'
$c.ClearFormats()

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = '@'
$c.Value = '-0.83%'
$c.ClearFormats()

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.8516'
$c.ClearFormats()

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = '@'
$c.Value = '-2.67%'
$c.ClearFormats()

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.9436'
$c.ClearFormats()

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = '@'
$c.Value = '-1.93%'
$c.ClearFormats()

$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = '@'
$c.Value = 'WazirX'
$c.ClearFormats()

$c = $ws.Cells.Item(9, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c.ClearFormats()

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.1389'
$c.ClearFormats()

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = '@'
$c.Value = '-1.75%'
$c.ClearFormats()

$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = '@'
$c.Value = 'LiechtensteinCryptoassetsExchange'
$c.ClearFormats()

$c = $ws.Cells.Item(10, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c.ClearFormats()

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.04970'
$c.ClearFormats()

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = '@'
$c.Value = '38.64%'
$c.ClearFormats()

$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = '@'
$c.Value = 'MandalaExchangeToken'
$c.ClearFormats()

$c = $ws.Cells.Item(11, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c.ClearFormats()

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.07087'
$c.ClearFormats()

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = '@'
$c.Value = '-1.31%'
$c.ClearFormats()

$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = '@'
$c.Value = 'BitrueCoin'
$c.ClearFormats()

$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c.ClearFormats()

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '0.03073'
$c.ClearFormats()

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = '@'
$c.Value = '-2.05%'
$c.ClearFormats()

$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = '@'
$c.Value = 'BitMartToken'
$c.ClearFormats()

$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c.ClearFormats()

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '0.09123'
$c.ClearFormats()

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = '@'
$c.Value = '-1.21%'
$c.ClearFormats()

$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = '@'
$c.Value = 'BitForexToken'
$c.ClearFormats()

$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c.ClearFormats()

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '0.001524'
$c.ClearFormats()

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '-1.38%'
$c.ClearFormats()

$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = '@'
$c.Value = 'One'
$c.ClearFormats()

$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c.ClearFormats()

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.0006040'
$c.ClearFormats()

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = '-0.13%'
$c.ClearFormats()

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '0.006024'
$c.ClearFormats()

$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = '0.34%'
$c.ClearFormats()

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '3.494'
$c.ClearFormats()

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = '0.26%'
$c.ClearFormats()

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '3.180'
$c.ClearFormats()

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = '-1.30%'
$c.ClearFormats()

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = '-0.67%'
$c.ClearFormats()

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '0.3053'
$c.ClearFormats()

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = '-2.93%'
$c.ClearFormats()

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = '-2.79%'
$c.ClearFormats()

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '3.951'
$c.ClearFormats()

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = '@'
$c.Value = '11.97%'
$c.ClearFormats()

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '0.04271'
$c.ClearFormats()

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = '1.13%'
$c.ClearFormats()

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = '-0.05%'
$c.ClearFormats()

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '0.004287'
$c.ClearFormats()

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = '-5.13%'
$c.ClearFormats()

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = '0.06%'
$c.ClearFormats()

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '0.0001523'
$c.ClearFormats()

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = '2.07%'
$c.ClearFormats()

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '0.03827'
$c.ClearFormats()

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = '@'
$c.Value = '-0.24%'
$c.ClearFormats()

$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = '@'
$c.Value = 'BKEXToken'
$c.ClearFormats()

$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c.ClearFormats()

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.1101'
$c.ClearFormats()

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = '@'
$c.Value = '-0.33%'
$c.ClearFormats()

$c = $ws.Cells.Item(42, 2)
$c.NumberFormat = '@'
$c.Value = 'KickToken'
$c.ClearFormats()

$c = $ws.Cells.Item(42, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$c.ClearFormats()

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.003897'
$c.ClearFormats()

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = '@'
$c.Value = '-33.79%'
$c.ClearFormats()

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.002430'
$c.ClearFormats()

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = '@'
$c.Value = '5.26%'
$c.ClearFormats()

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.01412'
$c.ClearFormats()

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = '@'
$c.Value = '34.58%'
$c.ClearFormats()

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.00005347'
$c.ClearFormats()

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = '@'
$c.Value = '-2.60%'
$c.ClearFormats()

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = '@'
$c.Value = '0.04%'
$c.ClearFormats()

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '0.06688'
$c.ClearFormats()

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = '@'
$c.Value = '-38.67%'
$c.ClearFormats()

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '0.2523'
$c.ClearFormats()

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = '@'
$c.Value = '11,611.46%'
$c.ClearFormats()

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '0.00002100'
$c.ClearFormats()

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = '@'
$c.Value = '0.04%'
$c.ClearFormats()

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '0.0002000'
$c.ClearFormats()

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = '@'
$c.Value = '0.04%'
$c.ClearFormats()
